# Rename the header row (row 1) from the generic "_old"/"_new" suffixes to
# the concrete format-version suffixes "_FV2304" (old/left side) and
# "_FV2310" (new/right side). Column K ("diff") is left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Segmentname_FV2304"
$ws.Range("B1").Value = "Segmentgruppe_FV2304"
$ws.Range("C1").Value = "Segment_FV2304"
$ws.Range("D1").Value = "Datenelement_FV2304"
$ws.Range("E1").Value = "Segment ID_FV2304"
$ws.Range("F1").Value = "Code_FV2304"
$ws.Range("G1").Value = "Qualifier_FV2304"
$ws.Range("H1").Value = "Beschreibung_FV2304"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2304"
$ws.Range("J1").Value = "Bedingung_FV2304"

$ws.Range("L1").Value = "Segmentname_FV2310"
$ws.Range("M1").Value = "Segmentgruppe_FV2310"
$ws.Range("N1").Value = "Segment_FV2310"
$ws.Range("O1").Value = "Datenelement_FV2310"
$ws.Range("P1").Value = "Segment ID_FV2310"
$ws.Range("Q1").Value = "Code_FV2310"
$ws.Range("R1").Value = "Qualifier_FV2310"
$ws.Range("S1").Value = "Beschreibung_FV2310"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2310"
$ws.Range("U1").Value = "Bedingung_FV2310"

# Turn the used range into a real table ("Table1") so the header row carries
# an AutoFilter + structured-table definition (xl/tables/table1.xml), with
# column names matching the renamed headers above.
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U60"), $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# Freeze the header row: selecting the first cell below the split before
# toggling FreezePanes yields a clean "frozen" pane (vs. a split/frozen
# mix) with the bottom-left pane anchored at A2.
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
